$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.4218206666666666
$ws.Cells.Item(2, 8).Value = 1.265462
$ws.Cells.Item(2, 9).Value = 0.204479520571209
$ws.Cells.Item(2, 10).Value = 0.204479520571209
$ws.Cells.Item(2, 13).Value = 1.815761
$ws.Cells.Item(2, 14).Value = 5.447283000000001
$ws.Cells.Item(2, 15).Value = 0.07007596730428067
$ws.Cells.Item(2, 16).Value = 0.07007596730428067
$ws.Cells.Item(2, 17).Value = 0.7659255155273333
$ws.Cells.Item(2, 18).Value = 6.893329639746
$ws.Cells.Item(2, 19).Value = 0.01432910019794303
$ws.Cells.Item(2, 20).Value = 0.01432910019794303
$ws.Cells.Item(3, 7).Value = 0.4218206666666666
$ws.Cells.Item(3, 8).Value = 1.265462
$ws.Cells.Item(3, 9).Value = 0.204479520571209
$ws.Cells.Item(3, 10).Value = 0.204479520571209
$ws.Cells.Item(3, 15).Value = 0.5079540516959071
$ws.Cells.Item(3, 16).Value = 0.5079540516959072
$ws.Cells.Item(3, 17).Value = 5.551902940134221
$ws.Cells.Item(3, 18).Value = 49.967126461208
$ws.Cells.Item(3, 19).Value = 0.1038662009629822
$ws.Cells.Item(3, 20).Value = 0.1038662009629822
$ws.Cells.Item(4, 7).Value = 0.4218206666666666
$ws.Cells.Item(4, 8).Value = 1.265462
$ws.Cells.Item(4, 9).Value = 0.204479520571209
$ws.Cells.Item(4, 10).Value = 0.204479520571209
$ws.Cells.Item(4, 13).Value = 9.711409333333334
$ws.Cells.Item(4, 14).Value = 29.134228
$ws.Cells.Item(4, 15).Value = 0.3747940411327002
$ws.Cells.Item(4, 16).Value = 0.3747940411327002
$ws.Cells.Item(4, 17).Value = 4.096473159259555
$ws.Cells.Item(4, 18).Value = 36.868258433336
$ws.Cells.Item(4, 19).Value = 0.07663770584376051
$ws.Cells.Item(4, 20).Value = 0.07663770584376051
$ws.Cells.Item(5, 7).Value = 0.4218206666666666
$ws.Cells.Item(5, 8).Value = 1.265462
$ws.Cells.Item(5, 9).Value = 0.204479520571209
$ws.Cells.Item(5, 10).Value = 0.204479520571209
$ws.Cells.Item(5, 13).Value = 1.222391
$ws.Cells.Item(5, 14).Value = 3.667173
$ws.Cells.Item(5, 15).Value = 0.04717593986711188
$ws.Cells.Item(5, 16).Value = 0.04717593986711189
$ws.Cells.Item(5, 17).Value = 0.5156297865473333
$ws.Cells.Item(5, 18).Value = 4.640668078926
$ws.Cells.Item(5, 19).Value = 0.009646513566523223
$ws.Cells.Item(5, 20).Value = 0.009646513566523223
$ws.Cells.Item(6, 9).Value = 0.3030684321645684
$ws.Cells.Item(6, 10).Value = 0.3030684321645683
$ws.Cells.Item(6, 13).Value = 1.815761
$ws.Cells.Item(6, 14).Value = 5.447283000000001
$ws.Cells.Item(6, 15).Value = 0.07007596730428067
$ws.Cells.Item(6, 16).Value = 0.07007596730428067
$ws.Cells.Item(6, 17).Value = 1.135213171946333
$ws.Cells.Item(6, 18).Value = 10.216918547517
$ws.Cells.Item(6, 19).Value = 0.0212378135433239
$ws.Cells.Item(6, 20).Value = 0.02123781354332389
$ws.Cells.Item(7, 9).Value = 0.3030684321645684
$ws.Cells.Item(7, 10).Value = 0.3030684321645683
$ws.Cells.Item(7, 15).Value = 0.5079540516959071
$ws.Cells.Item(7, 16).Value = 0.5079540516959072
$ws.Cells.Item(7, 19).Value = 0.1539448380591187
$ws.Cells.Item(7, 20).Value = 0.1539448380591187
$ws.Cells.Item(8, 9).Value = 0.3030684321645684
$ws.Cells.Item(8, 10).Value = 0.3030684321645683
$ws.Cells.Item(8, 13).Value = 9.711409333333334
$ws.Cells.Item(8, 14).Value = 29.134228
$ws.Cells.Item(8, 15).Value = 0.3747940411327002
$ws.Cells.Item(8, 16).Value = 0.3747940411327002
$ws.Cells.Item(8, 17).Value = 6.071569878063555
$ws.Cells.Item(8, 18).Value = 54.64412890257199
$ws.Cells.Item(8, 19).Value = 0.1135882424307102
$ws.Cells.Item(8, 20).Value = 0.1135882424307102
$ws.Cells.Item(9, 9).Value = 0.3030684321645684
$ws.Cells.Item(9, 10).Value = 0.3030684321645683
$ws.Cells.Item(9, 13).Value = 1.222391
$ws.Cells.Item(9, 14).Value = 3.667173
$ws.Cells.Item(9, 15).Value = 0.04717593986711188
$ws.Cells.Item(9, 16).Value = 0.04717593986711189
$ws.Cells.Item(9, 17).Value = 0.7642384457363333
$ws.Cells.Item(9, 18).Value = 6.878146011626999
$ws.Cells.Item(9, 19).Value = 0.01429753813141555
$ws.Cells.Item(9, 20).Value = 0.01429753813141555
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.469433
$ws.Cells.Item(10, 8).Value = 1.408299
$ws.Cells.Item(10, 9).Value = 0.2275598195290835
$ws.Cells.Item(10, 10).Value = 0.2275598195290835
$ws.Cells.Item(10, 13).Value = 1.815761
$ws.Cells.Item(10, 14).Value = 5.447283000000001
$ws.Cells.Item(10, 15).Value = 0.07007596730428067
$ws.Cells.Item(10, 16).Value = 0.07007596730428067
$ws.Cells.Item(10, 17).Value = 0.8523781335130001
$ws.Cells.Item(10, 18).Value = 7.671403201617
$ws.Cells.Item(10, 19).Value = 0.01594647447308807
$ws.Cells.Item(10, 20).Value = 0.01594647447308806
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.469433
$ws.Cells.Item(11, 8).Value = 1.408299
$ws.Cells.Item(11, 9).Value = 0.2275598195290835
$ws.Cells.Item(11, 10).Value = 0.2275598195290835
$ws.Cells.Item(11, 15).Value = 0.5079540516959071
$ws.Cells.Item(11, 16).Value = 0.5079540516959072
$ws.Cells.Item(11, 17).Value = 6.178565107990666
$ws.Cells.Item(11, 18).Value = 55.607085971916
$ws.Cells.Item(11, 19).Value = 0.1155899323329874
$ws.Cells.Item(11, 20).Value = 0.1155899323329874
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.469433
$ws.Cells.Item(12, 8).Value = 1.408299
$ws.Cells.Item(12, 9).Value = 0.2275598195290835
$ws.Cells.Item(12, 10).Value = 0.2275598195290835
$ws.Cells.Item(12, 13).Value = 9.711409333333334
$ws.Cells.Item(12, 14).Value = 29.134228
$ws.Cells.Item(12, 15).Value = 0.3747940411327002
$ws.Cells.Item(12, 16).Value = 0.3747940411327002
$ws.Cells.Item(12, 17).Value = 4.558856017574667
$ws.Cells.Item(12, 18).Value = 41.029704158172
$ws.Cells.Item(12, 19).Value = 0.08528806436073315
$ws.Cells.Item(12, 20).Value = 0.08528806436073313
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.469433
$ws.Cells.Item(13, 8).Value = 1.408299
$ws.Cells.Item(13, 9).Value = 0.2275598195290835
$ws.Cells.Item(13, 10).Value = 0.2275598195290835
$ws.Cells.Item(13, 13).Value = 1.222391
$ws.Cells.Item(13, 14).Value = 3.667173
$ws.Cells.Item(13, 15).Value = 0.04717593986711188
$ws.Cells.Item(13, 16).Value = 0.04717593986711189
$ws.Cells.Item(13, 17).Value = 0.5738306743029999
$ws.Cells.Item(13, 18).Value = 5.164476068727
$ws.Cells.Item(13, 19).Value = 0.01073534836227488
$ws.Cells.Item(13, 20).Value = 0.01073534836227488
$ws.Cells.Item(14, 7).Value = 0.546446
$ws.Cells.Item(14, 8).Value = 1.639338
$ws.Cells.Item(14, 9).Value = 0.2648922277351392
$ws.Cells.Item(14, 10).Value = 0.2648922277351391
$ws.Cells.Item(14, 13).Value = 1.815761
$ws.Cells.Item(14, 14).Value = 5.447283000000001
$ws.Cells.Item(14, 15).Value = 0.07007596730428067
$ws.Cells.Item(14, 16).Value = 0.07007596730428067
$ws.Cells.Item(14, 17).Value = 0.992215335406
$ws.Cells.Item(14, 18).Value = 8.929938018654001
$ws.Cells.Item(14, 19).Value = 0.01856257908992568
$ws.Cells.Item(14, 20).Value = 0.01856257908992568
$ws.Cells.Item(15, 7).Value = 0.546446
$ws.Cells.Item(15, 8).Value = 1.639338
$ws.Cells.Item(15, 9).Value = 0.2648922277351392
$ws.Cells.Item(15, 10).Value = 0.2648922277351391
$ws.Cells.Item(15, 15).Value = 0.5079540516959071
$ws.Cells.Item(15, 16).Value = 0.5079540516959072
$ws.Cells.Item(15, 17).Value = 7.192191833554666
$ws.Cells.Item(15, 18).Value = 64.729726501992
$ws.Cells.Item(15, 19).Value = 0.1345530803408189
$ws.Cells.Item(15, 20).Value = 0.1345530803408189
$ws.Cells.Item(16, 7).Value = 0.546446
$ws.Cells.Item(16, 8).Value = 1.639338
$ws.Cells.Item(16, 9).Value = 0.2648922277351392
$ws.Cells.Item(16, 10).Value = 0.2648922277351391
$ws.Cells.Item(16, 13).Value = 9.711409333333334
$ws.Cells.Item(16, 14).Value = 29.134228
$ws.Cells.Item(16, 15).Value = 0.3747940411327002
$ws.Cells.Item(16, 16).Value = 0.3747940411327002
$ws.Cells.Item(16, 17).Value = 5.306760784562667
$ws.Cells.Item(16, 18).Value = 47.760847061064
$ws.Cells.Item(16, 19).Value = 0.09928002849749631
$ws.Cells.Item(16, 20).Value = 0.0992800284974963
$ws.Cells.Item(17, 7).Value = 0.546446
$ws.Cells.Item(17, 8).Value = 1.639338
$ws.Cells.Item(17, 9).Value = 0.2648922277351392
$ws.Cells.Item(17, 10).Value = 0.2648922277351391
$ws.Cells.Item(17, 13).Value = 1.222391
$ws.Cells.Item(17, 14).Value = 3.667173
$ws.Cells.Item(17, 15).Value = 0.04717593986711188
$ws.Cells.Item(17, 16).Value = 0.04717593986711189
$ws.Cells.Item(17, 17).Value = 0.6679706723859999
$ws.Cells.Item(17, 18).Value = 6.011736051474
$ws.Cells.Item(17, 19).Value = 0.01249653980689823
$ws.Cells.Item(17, 20).Value = 0.01249653980689823
